# Add lesson 19 ("orderofoperation.xlsx") data to the "Completed Lessons" sheet.
# New rows are appended for Episode 20-25 and the two new lesson names that
# belong to rows 19-21 of column C are filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed Lessons")

# --- Fill in the lesson name for the already-existing "Episode 19" row ---
$ws.Range("C19").Value = "Relative vs Absolute Cell Reference in Formulas"

# --- Append new episode rows in column B ---
$ws.Range("B21").Value = "Episode 20"
$ws.Range("B22").Value = "Episode 21"
$ws.Range("B23").Value = "Episode 22"
$ws.Range("B24").Value = "Episode 23"
$ws.Range("B25").Value = "Episode 24"
$ws.Range("B26").Value = "Episode 25"

# --- Fill in the lesson names in column C for rows 20 and 21 ---
$ws.Range("C20").Value = "Understanding the Order of Operation"
$ws.Range("C21").Value = "The Structure of an Excel Function"

# --- Match row heights used throughout the rest of the sheet ---
$ws.Rows.Item(21).RowHeight = 23.4
$ws.Rows.Item(22).RowHeight = 23.4
$ws.Rows.Item(23).RowHeight = 23.4
$ws.Rows.Item(24).RowHeight = 23.4
$ws.Rows.Item(25).RowHeight = 23.4
$ws.Rows.Item(26).RowHeight = 23.4

# --- Update the active selection to reflect where the user ended up ---
$ws.Activate()
$ws.Range("C22").Select()
